$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Insert new "Brave" row under the Browser category (old row 24 -> 25) ---
$ws1.Rows.Item(24).Insert()

$ws1.Cells.Item(24,1).Value = "Browser"
$ws1.Cells.Item(24,2).Value = "Brave"
$ws1.Cells.Item(24,3).Value = "choco install brave --pre "
$ws1.Cells.Item(24,4).Value = "winget install --id=Brave.Brave  -e"
$ws1.Cells.Item(24,5).Value = "scoop bucket add extras `nscoop install extras/brave"
$ws1.Cells.Item(24,6).Value = "sudo apt install curl`nsudo curl -fsSLo /usr/share/keyrings/brave-browser-archive-keyring.gpg https://brave-browser-apt-release.s3.brave.com/brave-browser-archive-keyring.gpg`nsudo curl -fsSLo /etc/apt/sources.list.d/brave-browser-release.sources https://brave-browser-apt-release.s3.brave.com/brave-browser.sources`nsudo apt update`nsudo apt install brave-browser"
$ws1.Cells.Item(24,7).Value = "sudo dnf install dnf-plugins-core  sudo dnf config-manager addrepo --from-repofile=https://brave-browser-rpm-release.s3.brave.com/brave-browser.repo  sudo dnf install brave-browser"
$ws1.Cells.Item(24,8).Value = "sudo pacman -Sy brave-browser"
$ws1.Cells.Item(24,9).Value = "brew install --cask brave-browser"
$ws1.Rows.Item(24).RowHeight = 686.05

# --- Kubernetes row (shifted to row 26): dnf cell repurposed with brave repo command ---
$ws1.Cells.Item(26,7).Value = "sudo dnf config-manager addrepo --from-repofile=https://brave-browser-rpm-release.s3.brave.com/brave-browser.repo"
$ws1.Rows.Item(26).RowHeight = 68.65

# --- Ansible row (shifted to row 27): content unchanged, keep its original height ---
$ws1.Rows.Item(27).RowHeight = 41.75

# --- Nmap row (shifted to row 28): dnf cell repurposed to install brave-browser ---
$ws1.Cells.Item(28,7).Value = "sudo dnf install brave-browser"
$ws1.Rows.Item(28).RowHeight = 68.65

# --- Wireshark row (shifted to row 29): content unchanged, keep its original height ---
$ws1.Rows.Item(29).RowHeight = 55.2

# --- New "Utilties" category rows ---
# Row 30
$ws1.Cells.Item(30,1).Value = "Utilties"
$ws1.Cells.Item(30,2).Value = "Vlc"
$ws1.Cells.Item(30,3).Value = "choco install vlc"
$ws1.Cells.Item(30,4).Value = "winget install --id=VideoLAN.VLC  -e"
$ws1.Cells.Item(30,5).Value = "scoop bucket add extras`nScoop install extras/vlc"
$ws1.Cells.Item(30,6).Value = "sudo apt install vlc"
$ws1.Cells.Item(30,7).Value = "sudo dnf install vlc"
$ws1.Cells.Item(30,8).Value = "sudo pacman -S vlc"
$ws1.Cells.Item(30,9).Value = "brew install --cask vlc"
$ws1.Cells.Item(30,5).WrapText = $true
$ws1.Cells.Item(30,7).WrapText = $true
$ws1.Cells.Item(30,8).WrapText = $true
$ws1.Rows.Item(30).RowHeight = 41.75

# Row 31
$ws1.Cells.Item(31,1).Value = "Utilties"
$ws1.Cells.Item(31,2).Value = "obs studio"
$ws1.Cells.Item(31,3).Value = "choco install obs-studio"
$ws1.Cells.Item(31,4).Value = "winget install --id=OBSProject.OBSStudio  -e"
$ws1.Cells.Item(31,5).Value = "scoop bucket add extras`nscoop install extras/obs-studio"
$ws1.Cells.Item(31,6).Value = "t: sudo apt install obs-studio"
$ws1.Cells.Item(31,7).Value = "sudo dnf install obs-studio"
$ws1.Cells.Item(31,8).Value = "sudo pacman -S obs-studio"
$ws1.Cells.Item(31,5).WrapText = $true
$ws1.Cells.Item(31,7).WrapText = $true
$ws1.Cells.Item(31,8).WrapText = $true
$ws1.Rows.Item(31).RowHeight = 55.2

# Row 32
$ws1.Cells.Item(32,1).Value = "Utilties"
$ws1.Cells.Item(32,2).Value = "notion"
$ws1.Cells.Item(32,3).Value = "choco install notion"
$ws1.Cells.Item(32,4).Value = "winget install --id=Notion.Notion  -e"
$ws1.Cells.Item(32,5).Value = "scoop bucket add extras`nscoop install extras/notion"
$ws1.Cells.Item(32,9).Value = "brew install --cask notion"
$ws1.Cells.Item(32,5).WrapText = $true
$ws1.Rows.Item(32).RowHeight = 28.35

# Row 33
$ws1.Cells.Item(33,1).Value = "Utilties"
$ws1.Cells.Item(33,2).Value = "ollama"
$ws1.Cells.Item(33,3).Value = "choco install ollama"
$ws1.Cells.Item(33,4).Value = "winget install --id=Ollama.Ollama  -e"
$ws1.Cells.Item(33,5).Value = "scoop bucket add main`nscoop install main/ollama"
$ws1.Cells.Item(33,9).Value = "brew install ollama"
$ws1.Cells.Item(33,5).WrapText = $true
$ws1.Rows.Item(33).RowHeight = 28.35

# --- View state changes ---
$excel.ActiveWindow.Zoom = 88
$ws1.Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("I24").Select()
